$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.176.53"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "2.274.40"
$ws.Range("E3").Value = "  -2.60%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "298.20"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "94.86"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.85%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -3.73%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.491"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.85%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "33.16"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.73%  "
$ws.Range("E11").Value = "  -1.09%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "48.18"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -8.21%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "15.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.63"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "2.626.17"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").Value = "2.258.33"
$ws.Range("E17").Value = "  -3.30%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.777"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").Value = "42.156.97"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "0.0₃0889"
$ws.Range("E20").Value = "  -2.33%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.39"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.47%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.90%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "66.60"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.93%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "232.70"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.72%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.68%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -4.50%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "23.85"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.12%  "
$ws.Range("E29").Value = "  -1.17%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "166.77"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.82%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "33.67"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.55%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.04"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("E33").Value = "  -0.11%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("E36").Value = "  -5.01%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0690"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.23%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "16.04"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -8.43%  "
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("E42").Value = "  -7.11%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.41"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.37%  "
$ws.Range("D44").Value = "1.958.36"
$ws.Range("E44").Value = "  -3.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0278"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.71%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.46"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -6.98%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.57"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("E48").Value = "  -4.74%  "
$ws.Range("D49").Value = "2.497.67"
$ws.Range("E49").Value = "  -2.18%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "51.97"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -7.41%  "
$ws.Range("E51").Value = "  -5.19%  "
